$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Routine_Care/Nursing for Arterial and Central Venous Lines.pdf', '2011-03-06')
    ,@(3, 'Routine_Care/VTE_Prevention/TED Stocking Sizing.pdf', '2011-06-06')
    ,@(4, 'Routine_Care/Faecal  incontinence skin care.pdf', '2011-07-06')
    ,@(5, 'Breathing(Respiratory)/Equipment/IPPB using an ICU Ventilator.pdf', '2011-12-06')
    ,@(6, 'Drugs/heparin_critical_care_only.pdf', '2014-09-06')
    ,@(7, 'Routine_Care/Invasive Flush Systems.pdf', '2014-10-06')
    ,@(8, 'GI_Liver_and_Transplant/Pancreatic Irrigation.pdf', '2014-11-06')
    ,@(9, 'Neurological/Thiopentone levels.pdf', '2014-12-06')
    ,@(10, 'Infection_and_sepsis/Ebola/Ebola.pdf', '2015-03-06')
    ,@(11, 'GI_Liver_and_Transplant/Nasal bridle.pdf', '2016-01-06')
    ,@(12, 'Breathing(Respiratory)/Equipment/AMBU AScope.pdf', '2016-04-06')
    ,@(13, 'Cardiovascular/EZ-IO Intraosseus Access Device_pub_em.pdf', '2016-08-06')
    ,@(14, 'Routine_Care/Central venous catheter removal.pdf', '2017-01-06')
    ,@(15, 'Routine_Care/Tracheostomy_nursing_care.pdf', '2017-08-06')
    ,@(16, 'ECLS/Extra Corporeal Carbon Dioxide Removal.pdf', '2017-10-06')
    ,@(17, 'Airway/Critical care extubation checklist.pdf', '2017-10-06')
    ,@(18, 'Airway/Tracheostomy_Laryngectomy/Hospital_in-patients_with_a_Tracheostomy.pdf', '2017-11-06')
    ,@(19, 'Trauma and Burns/Mangement of burns.pdf', '2018-05-06')
    ,@(20, 'Drugs/diazepam_diazemuls.pdf', '2018-07-06')
    ,@(21, 'End_of_life_care/Reasons to report a death to PF.pdf', '2019-05-06')
    ,@(22, 'Breathing(Respiratory)/salbutamol and ipratroprium MDI.pdf', '2019-05-06')
    ,@(23, 'Neurological/SOP -  Femoral site care.pdf', '2019-06-06')
    ,@(24, 'Airway/Tracheostomy_Laryngectomy/Tracheostomy suctioning cleaning guideline.pdf', '2019-06-06')
    ,@(25, 'Airway/Tracheostomy_Laryngectomy/Tracheostomy change in Critical Care.pdf', '2019-06-06')
    ,@(26, 'Drugs/sodium_bicarbonate.pdf', '2019-08-06')
    ,@(27, 'Post_op_care/Anticoagulation antiplatelet agents and epidural analgesia.pdf', '2019-09-06')
    ,@(28, 'Post_op_care/Epidural top-up.pdf', '2020-01-06')
    ,@(29, 'Drugs/heparin for Haemofiltration.pdf', '2020-03-06')
    ,@(30, 'Covid-19/WGH/CoVid intubation checklist WGH.pdf', '2020-03-06')
    ,@(31, 'Covid-19/SJH/SJH COVID19 ITU Intubation Action Card.pdf', '2020-03-06')
    ,@(32, 'Diabetes_and_Glucose/Hyperosmolar Hyperglycaemic State.pdf', '2020-03-06')
    ,@(33, 'Covid-19/SJH/SJH COVID19 ED Intubation Action Card.pdf', '2020-03-06')
    ,@(34, 'Airway/Emergency intubation checklist_em_pub.pdf', '2020-03-06')
    ,@(35, 'Drugs/fentanyl.pdf', '2020-04-06')
    ,@(36, 'Airway/Tracheostomy_Laryngectomy/Tracheostomy guideline.pdf', '2020-05-06')
    ,@(37, 'Covid-19/WGH/WGH_CT_Transfer_May.pdf', '2020-07-06')
    ,@(38, 'Cardiovascular/GJNH Acute Heart Failure Referral Form.pdf', '2020-08-06')
    ,@(39, 'Organ_donation/Donation after circulatory death.pdf', '2020-11-06')
    ,@(40, 'Airway/Percutaneous tracheostomy checklist.pdf', '2021-02-06')
    ,@(41, 'Infection_and_sepsis/SOP Ultrasound Cleaning.pdf', '2021-05-06')
    ,@(42, 'Delirium/Managing a Potentially Violent Patient.pdf', '2021-05-06')
    ,@(43, 'Delirium/Risk assessment posi mit.pdf', '2021-05-06')
    ,@(44, 'Neurological/Sub arachnoid haemorrhage management.pdf', '2021-06-06')
    ,@(45, 'Airway/Anticipated difficult airway tool.pdf', '2021-06-06')
    ,@(46, 'Airway/McGrath Mac.pdf', '2021-06-06')
    ,@(47, 'Delirium/Drugs Causing Delirium and Agitiation.pdf', '2021-06-06')
    ,@(48, 'Airway/Tracheostomy_Laryngectomy/Tracheostomy safety box contents.pdf', '2021-06-06')
    ,@(49, 'Breathing(Respiratory)/HFNO.pdf', '2021-06-06')
    ,@(50, 'GI_Liver_and_Transplant/Treatment of constipation.pdf', '2021-06-06')
    ,@(51, 'GI_Liver_and_Transplant/Abdominal pressure measurement.pdf', '2021-06-06')
    ,@(52, 'Drugs/ketamine_in_asthma.pdf', '2021-06-06')
    ,@(53, 'End_of_life_care/Documentation following death.pdf', '2021-09-06')
    ,@(54, 'Drugs/zanamivir.pdf', '2021-12-06')
    ,@(55, 'Routine_Care/bBraun Spaceplus Failure EMERGENCY ACTION CARD_em.pdf', '2022-01-06')
    ,@(56, 'Drugs/insulin.pdf', '2022-03-06')
    ,@(57, 'Breathing(Respiratory)/Equipment/HFNO Set Up.pdf', '2022-03-06')
    ,@(58, 'Breathing(Respiratory)/Inhaled Nitrous Oxide.pdf', '2022-04-06')
    ,@(59, 'Breathing(Respiratory)/Equipment/APRV.pdf', '2022-05-06')
    ,@(60, 'Cardiovascular/Steroids for Septic Shock.pdf', '2022-05-06')
    ,@(61, 'Neurological/SOP for review of Neurosurgical patients in ITU by neurosurgical team.pdf', '2022-06-06')
    ,@(62, 'Post_op_care/Epidural Haematoma.pdf', '2022-06-06')
    ,@(63, 'Breathing(Respiratory)/Equipment/T piece Y piece.pdf', '2022-06-06')
    ,@(64, 'Policies_and_admin/General Critical Care Interaction with HEPMA_pub.pdf', '2022-07-06')
    ,@(65, 'Drugs/midazolam and thiopental levels.pdf', '2022-08-06')
    ,@(66, 'Breathing(Respiratory)/Equipment/HFNO through ventilator.pdf', '2022-10-06')
    ,@(67, 'Post_op_care/Adult Scoliosis Spinal Surgery Post-Op Care.pdf', '2022-11-06')
    ,@(68, 'Routine_Care/VTE_Prevention/Dalteparin_thromboprophylaxis.pdf', '2022-11-06')
    ,@(69, 'Post_op_care/Post op care pharyngo-laryngo-oesphagectomy PLOG.pdf', '2022-12-06')
    ,@(70, 'GI_Liver_and_Transplant/Nasogastric feeding protocol.pdf', '2023-01-06')
    ,@(71, 'Diabetes_and_Glucose/Intravenous Insulin Therapy (not for DKA or HHS).pdf', '2023-02-06')
    ,@(72, 'Drugs/Antibiotic doses in CVVHD.pdf', '2023-02-06')
    ,@(73, 'GI_Liver_and_Transplant/Jejunostomy feeding protocol.pdf', '2023-04-06')
    ,@(74, 'ECLS/RIE ECLS Anti Xa Protocol.pdf', '2023-04-06')
    ,@(75, 'GI_Liver_and_Transplant/Nasojejunal feeding protocol.pdf', '2023-04-06')
    ,@(76, 'Drugs/vasopressin_sepsis.pdf', '2023-05-06')
    ,@(77, 'Drugs/vasopressin organ donation.pdf', '2023-05-06')
    ,@(78, 'Infection_and_sepsis/Winter Infections Stepdown Guidance.pdf', '2023-05-06')
    ,@(79, 'Transfer/ACCP Transfers.pdf', '2023-06-06')
    ,@(80, 'Covid-19/videos/Donning and Doffing Video.pdf', '2023-06-06')
    ,@(81, 'Breathing(Respiratory)/Equipment/Ventilators Circuits Filters and Closed Suction - Set up and Maintenance.pdf', '2023-07-06')
    ,@(82, 'Breathing(Respiratory)/Equipment/Bipap V60.pdf', '2023-07-06')
    ,@(83, 'Infection_and_sepsis/Infection indications for IVIG.pdf', '2023-07-06')
    ,@(84, 'Breathing(Respiratory)/CPAP.pdf', '2023-07-06')
    ,@(85, 'Drugs/piperacillin_tazobactam extended_infusion.pdf', '2023-07-06')
    ,@(86, 'Procedures/CVC Guidance/Securing CVCs.pdf', '2023-08-06')
    ,@(87, 'Covid-19/Covid 19 Death Certification Guideline.pdf', '2023-08-06')
    ,@(88, 'Routine_Care/Video Communication.pdf', '2023-09-06')
    ,@(89, 'Neurological/Treatment of status epilepticus.pdf', '2023-09-06')
    ,@(90, 'Drugs/isoprenaline.pdf', '2023-10-06')
    ,@(91, 'Cardiovascular/Cardiogenic Shock.pdf', '2023-10-06')
    ,@(92, 'Haematology_CAR-T/Haem_ICU_transfer.pdf', '2024-01-06')
    ,@(93, 'Cardiovascular/Management of hypertension within Critical Care.pdf', '2024-02-06')
    ,@(94, 'Drugs/aminophylline.pdf', '2024-02-06')
    ,@(95, 'Haematology_CAR-T/ICANS.pdf', '2024-03-06')
    ,@(96, 'Drugs/pancuronium.pdf', '2024-03-06')
    ,@(97, 'Drugs/rocuronium.pdf', '2024-03-06')
    ,@(98, 'Haematology_CAR-T/CRS.pdf', '2024-03-06')
    ,@(99, 'Drugs/phenytoin.pdf', '2024-03-06')
    ,@(100, 'Policies_and_admin/General Critical Care SOP_pub.pdf', '2024-04-06')
    ,@(101, 'Drugs/milrinone.pdf', '2024-04-06')
    ,@(102, 'Neurological/Critical Care MRI Procedure_pub.pdf', '2024-05-06')
    ,@(103, 'End_of_life_care/Palliative extubation & withdrawal of invasive ventilatory support nursing checklist.pdf', '2024-05-06')
    ,@(104, 'End_of_life_care/CMO & NRS Guidance for Doctors completing MCCD.pdf', '2024-05-06')
    ,@(105, 'Infection_and_sepsis/Trip Out of Unit infection guidance.pdf', '2024-05-06')
    ,@(106, 'Drugs/clonidine.pdf', '2024-05-06')
    ,@(107, 'Ethics_and_Law/Care at the End of Life (FICM).pdf', '2024-05-06')
    ,@(108, 'GI_Liver_and_Transplant/ICU - Upper GI bleeding (Endoscopy guideline).pdf', '2024-05-06')
    ,@(109, 'Neurological/Management of traumatic brain injury.pdf', '2024-05-06')
    ,@(110, 'Ethics_and_Law/DNACPR policy for Scotland.pdf', '2024-05-06')
    ,@(111, 'Covid-19/COVID 19 ICM guidance basic goals_June_2022.pdf', '2024-05-06')
    ,@(112, 'Organ_donation/Organ Retrieval SOP.pdf', '2024-05-06')
    ,@(113, 'Airway/Cook Staged Extubation Set.pdf', '2024-06-02')
    ,@(114, 'Drugs/noradrenaline (central).pdf', '2024-06-06')
    ,@(115, 'Post_op_care/Epidural hypotension.pdf', '2024-06-26')
    ,@(116, 'Drugs/dexmedetomidine.pdf', '2024-07-06')
    ,@(117, 'GI_Liver_and_Transplant/Fulminant Liver Failure.pdf', '2024-07-06')
    ,@(118, 'GI_Liver_and_Transplant/Confirmation of Nasogastric Tube Position.pdf', '2024-07-06')
    ,@(119, 'Breathing(Respiratory)/Equipment/Passy Muir Valve.pdf', '2024-07-06')
    ,@(120, 'Drugs/glyceryl_trinitrate.pdf', '2024-07-06')
    ,@(121, 'Infection_and_sepsis/Antifungal guidance in critical care.pdf', '2024-07-25')
    ,@(122, 'Cardiovascular/Intra Aortic Balloon Pump Bedside Checks_pub.pdf', '2024-08-06')
    ,@(123, 'Cardiovascular/Intra Aortic Balloon Pump Guideline_pub.pdf', '2024-08-06')
    ,@(124, 'Transfer/Transfer Outdoors to Garden Guideline.pdf', '2024-08-06')
    ,@(125, 'Breathing(Respiratory)/ARDS Strategy.pdf', '2024-08-15')
    ,@(126, 'Drugs/adrenaline.pdf', '2024-10-06')
    ,@(127, 'Drugs/dobutamine.pdf', '2024-10-06')
    ,@(128, 'Drugs/hydralazine.pdf', '2024-10-24')
    ,@(129, 'Post_op_care/Major OMFS Free Flap.pdf', '2024-11-20')
    ,@(130, 'Drugs/Alteplase for massive PE.pdf', '2024-11-24')
    ,@(131, 'Drugs/alfentanil.pdf', '2024-11-24')
    ,@(132, 'Drugs/magnesium.pdf', '2024-12-24')
    ,@(133, 'Drugs/all IV drug infusion information.pdf', '2025-01-06')
    ,@(134, 'Drugs/vancomycin.pdf', '2025-01-06')
    ,@(135, 'Drugs/neostigmine.pdf', '2025-01-06')
    ,@(136, 'Drugs/labetalol.pdf', '2025-02-06')
    ,@(137, 'Infection_and_sepsis/Initial investigation and management in unidentified Infections.pdf', '2025-02-06')
    ,@(138, 'Neurological/Intrathecal policy RIE.pdf', '2025-02-06')
    ,@(139, 'Cardiovascular/Management of Acute Type B Aortic Dissection Guideline.pdf', '2025-03-06')
    ,@(140, 'Drugs/midazolam.pdf', '2025-03-06')
    ,@(141, 'Procedures/CVC Guidance/CVC NHL  April 2023.pdf', '2025-04-06')
    ,@(142, 'Drugs/salbutamol.pdf', '2025-04-06')
    ,@(143, 'Drugs/nimodipine.pdf', '2025-04-06')
    ,@(144, 'Drugs/potassium.pdf', '2025-04-06')
    ,@(145, 'Drugs/nicardipine.pdf', '2025-05-06')
    ,@(146, 'Drugs/phenobarbitone.pdf', '2025-05-06')
    ,@(147, 'Routine_Care/ICU Eye Care Guideline.pdf', '2025-05-06')
    ,@(148, 'Drugs/amiodarone.pdf', '2025-05-06')
    ,@(149, 'Procedures/Arterial Line insertion for ACCPs.pdf', '2025-05-06')
    ,@(150, 'Drugs/phenylephrine.pdf', '2025-06-06')
    ,@(151, 'Breathing(Respiratory)/Manual Ventilation and MHI.pdf', '2025-06-06')
    ,@(152, 'Drugs/morphine.pdf', '2025-06-06')
    ,@(153, 'Drugs/noradrenaline (peripheral).pdf', '2025-06-06')
    ,@(154, 'Drugs/Epoprostenol.pdf', '2025-06-06')
    ,@(155, 'Neurological/Ventriculitis Guideline.pdf', '2025-06-06')
    ,@(156, 'Drugs/calcium.pdf', '2025-07-06')
    ,@(157, 'Cardiovascular/Cardiac Output Monitoring_pub .pdf', '2025-07-06')
    ,@(158, 'Cardiovascular/Pulmonary_Embolism_and_DVT/Catheter directed thrombolysis of iliofemoral DVT alteplase_pub.pdf', '2025-07-06')
    ,@(159, 'Drugs/atracurium.pdf', '2025-08-06')
    ,@(160, 'Drugs/dalteparin_thromboprophylaxis.pdf', '2025-08-06')
    ,@(161, 'Airway/Tracheostomy_Laryngectomy/Decannulation Guidline.pdf', '2025-08-06')
    ,@(162, 'Drugs/Vancomycin Continuous Infusion Fluid Restricted.pdf', '2025-08-06')
    ,@(163, 'Policies_and_admin/Pet Visitation.pdf', '2025-10-06')
    ,@(164, 'Drugs/valproate.pdf', '2025-10-06')
    ,@(165, 'Procedures/Inadvertent Catheter Placement Guideline.pdf', '2025-10-06')
    ,@(166, 'Drugs/Phosphate.pdf', '2025-11-06')
    ,@(167, 'Drugs/stress ulcer prophylaxis.pdf', '2025-11-06')
    ,@(168, 'Drugs/ketamine_for_status epilepticus.pdf', '2025-11-06')
    ,@(169, 'GI_Liver_and_Transplant/Prokinetics in ICU.pdf', '2026-01-06')
    ,@(170, 'Drugs/Thiopentone.pdf', '2026-01-06')
    ,@(171, 'Drugs/Octreotide.pdf', '2026-01-06')
    ,@(172, 'Breathing(Respiratory)/Proning Guideline.pdf', '2026-01-06')
    ,@(173, 'Procedures/ACCPs acquiring initial CVC competencies.pdf', '2026-03-06')
    ,@(174, 'Procedures/ACCP CVC placement following completion of initial competencies.pdf', '2026-03-06')
    ,@(175, 'Post_op_care/Prevention and treatment of paraplegia after major aortic procedures.pdf', '2026-03-06')
    ,@(176, 'Delirium/Violence and Agression.pdf', '2026-05-06')
    ,@(177, 'Post_op_care/Care of the Transgender Patient.pdf', '2026-07-06')
    ,@(178, 'GI_Liver_and_Transplant/Plasma exchange in Acute Liver Failure.pdf', '2026-11-06')
    ,@(179, 'Breathing(Respiratory)/Equipment/NIV through Drager Vent Set up in Critical Care.pdf', '2027-01-06')
    ,@(180, 'End_of_life_care/Guideline following Sudden Cardiac Death where death occurs in ICU.pdf', '2027-01-06')
    ,@(181, 'Transfer/Transfer Guidelines.pdf', '2027-02-06')
    ,@(182, 'Breathing(Respiratory)/Equipment/NIV through Nihon Kohden  Setup.pdf', '2028-02-06')
    ,@(183, 'Infection_and_sepsis/Influenza in Critical Care.pdf', '2028-05-06')
    ,@(184, 'Policies_and_admin/Anticipated Post op flow surgical patients.pdf', '2028-07-06')
    ,@(185, 'Policies_and_admin/Discharge Home from Critical Care.pdf', '2028-11-06')
    ,@(186, 'Policies_and_admin/Repatriaiton Checklist for Critical Care.pdf', '2030-11-06')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = "'" + $row[2]
}

